$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: refresh header/data (name/last name/age/salary -> Name/age/City/salary)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A1").Value = "Name"
$ws1.Range("B1").Value = "age"
$ws1.Range("C1").Value = "City"
$ws1.Range("D1").Value = "salary"

$ws1.Range("A2").Value = "Hamid"
$ws1.Range("B2").Value = 30
$ws1.Range("C2").Value = "San Ramon"
$ws1.Range("D2").Value = 120000

$ws1.Range("A3").Value = "Sam"
$ws1.Range("B3").Value = 60
$ws1.Range("C3").Value = "Houston"
$ws1.Range("D3").Value = 125000

$ws1.Range("A4").Value = "Fayed"
$ws1.Range("B4").Value = 45
$ws1.Range("C4").Value = "Alexandria"
$ws1.Range("D4").Value = 129000

# B4 used to carry a special number-format style (340009) -- drop it back to Normal.
$ws1.Range("B4").Style = "Normal"

# View: zoomed to 181%, whole grid selected, no single cell highlighted.
$ws1.Range("A1:D4").Select()
$excel.ActiveWindow.Zoom = 181

# ---------------------------------------------------------------------------
# Sheet2 (new): Username/Password table with hyperlinked emails
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Useraname"
$ws2.Range("B1").Value = "Password"

# Row 2: valid email - hyperlink display matches the cell text, so no
# explicit TextToDisplay is supplied (keeps the serialized XML free of a
# redundant display= attribute).
$ws2.Range("A2").Value = "asghar@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:asghar@gmail.com") | Out-Null
$ws2.Range("B2").Value = "123SKBDSSD+_#_$"

# Row 3: missing ".com"
$ws2.Range("A3").Value = "asghar@gmail"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:asghar@gmail") | Out-Null
$ws2.Range("B3").Value = "123SKBDSSD+_#_$"

# Row 4: trailing dot
$ws2.Range("A4").Value = "asghar@gmail."
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:asghar@gmail.") | Out-Null
$ws2.Range("B4").Value = "123SKBDSSD+_#_$"

# Rows 5-7: missing "@" in the cell text; the underlying hyperlink still
# targets "asghar@gmail." (added first with that text, then the cell text
# is overwritten), which leaves the stale display="asghar@gmail." seen in
# the source file.
$ws2.Hyperlinks.Add($ws2.Range("A5"), "mailto:asghar@gmail.", "", "", "asghar@gmail.") | Out-Null
$ws2.Range("A5").Value = "asghargmail.com"
$ws2.Range("B5").Value = "123SKBDSSD+_#_$"

$ws2.Hyperlinks.Add($ws2.Range("A6"), "mailto:asghar@gmail.", "", "", "asghar@gmail.") | Out-Null
$ws2.Range("A6").Value = "asghargmail.com"
$ws2.Range("B6").Value = "123SKBDSSD"

$ws2.Hyperlinks.Add($ws2.Range("A7"), "mailto:asghar@gmail.", "", "", "asghar@gmail.") | Out-Null
$ws2.Range("A7").Value = "asghargmail.com"
$ws2.Range("B7").Value = "123SKBDSSD"

# Column widths (best-fit-ish).
$ws2.Columns.Item(1).ColumnWidth = 16.5
$ws2.Columns.Item(2).ColumnWidth = 20.83203125

# View: zoomed to 132%, active cell parked below the data (B21), Sheet2 is
# the tab that's on top when the workbook is reopened.
$ws2.Range("B21").Select()
$excel.ActiveWindow.Zoom = 132
$ws2.Activate()
